$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 57 (shifts old rows 57-98 down to 59-100).
$ws.Rows("57:58").Insert()

# Populate the first new row (57) with the new weekly entry.
$ws.Range("A57").Value = 8
$ws.Range("B57").Value = "Terminal La Palmera de La Serena"
$ws.Range("C57").Value = "Coquimbo"
$ws.Range("D57").Value = 45079
$ws.Range("E57").Value = 4
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100104
$ws.Range("H57").Value = "Frutos de pepita"
$ws.Range("I57").Value = 100104003
$ws.Range("J57").Value = "Membrillo"
$ws.Range("K57").Value = "Champion"
$ws.Range("L57").Value = "Primera"
$ws.Range("M57").Value = 16
$ws.Range("N57").Value = 220000
$ws.Range("O57").Value = 230000
$ws.Range("P57").Value = 225000
$ws.Range("Q57").Value = "$/bins (450 kilos)"
$ws.Range("R57").Value = "Región de O'Higgins"
$ws.Range("S57").Value = 500
$ws.Range("T57").Value = 450

# Populate the second new row (58) with the new weekly entry.
$ws.Range("A58").Value = 8
$ws.Range("B58").Value = "Terminal La Palmera de La Serena"
$ws.Range("C58").Value = "Coquimbo"
$ws.Range("D58").Value = 45079
$ws.Range("E58").Value = 4
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100104
$ws.Range("H58").Value = "Frutos de pepita"
$ws.Range("I58").Value = 100104003
$ws.Range("J58").Value = "Membrillo"
$ws.Range("K58").Value = "Champion"
$ws.Range("L58").Value = "Segunda"
$ws.Range("M58").Value = 16
$ws.Range("N58").Value = 190000
$ws.Range("O58").Value = 200000
$ws.Range("P58").Value = 195000
$ws.Range("Q58").Value = "$/bins (450 kilos)"
$ws.Range("R58").Value = "Región de O'Higgins"
$ws.Range("S58").Value = 433
$ws.Range("T58").Value = 450

Write-Output $ws.UsedRange.Address()
